$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy column formatting from BK (63) to the new BL (64) column so the new
#    cells inherit the correct style indices (header row + the 257 data/total
#    rows). Rows below 259 are left untouched on purpose - they get no new
#    cell content, only their `spans` bookkeeping changes (handled below via
#    UsedRange).
$src = $ws.Range("BK2:BK259")
$dst = $ws.Range("BL2:BL259")
$src.Copy()
$dst.PasteSpecial(-4122)

# 2. New column header (row 2).
$ws.Cells.Item(2, 64).Value = "Tests Through June 22"

# 3. New column data, rows 3-259 (one value per county + the statewide total
#    in row 259).
$blValues = @(5058, 95, 3662, 342, 94, 180, 1424, 1418, 107, 425, 2285, 73, 1734, 22592, 63752, 2044, 2, 998, 6910, 14138, 12295, 943, 56, 111, 1061, 703, 1752, 1145, 1013, 315, 13833, 356, 152, 778, 318, 1441, 1392, 289, 325, 19, 36, 190, 39489, 159, 1183, 5302, 681, 93, 1607, 6734, 37, 157, 109, 159, 38, 10, 152018, 444, 424, 101, 35563, 520, 26, 187, 190, 386, 419, 4712, 31, 9005, 44256, 1516, 847, 1432, 1466, 32, 279, 96, 89616, 262, 680, 748, 275, 38251, 3, 1005, 72, 177, 1091, 862, 7378, 3027, 1106, 3736, 436, 161, 784, 55, 60, 903, 257823, 4654, 89, 286, 10923, 182, 1959, 23699, 1267, 94, 1227, 721, 1588, 399, 58, 3903, 637, 6, 106, 884, 869, 11, 15052, 132, 2211, 8010, 605, 322, 4649, 1296, 3, 3, 2293, 130, 4, 23, 2343, 210, 1361, 216, 1536, 529, 677, 369, 569, 3438, 923, 120, 391, 1149, 1, 13101, 475, 92, 13407, 11, 582, 105, 110, 353, 1667, 3408, 1135, 31, 5273, 1630, 108, 226, 326, 26283, 2904, 238, 27, 2954, 1684, 315, 161, 6378, 145, 113, 1545, 674, 643, 5367, 326, 1136, 1915, 14204, 133, 116, 5019, 115, 139, 418, 235, 213, 7, 701, 4882, 335, 1291, 517, 298, 691, 1210, 197, 106, 103, 44, 1477, 35, 7212, 290, 2639, 113, 27, 9, 64, 254, 89133, 8087, 6, 271, 43, 2850, 2731, 56384, 692, 662, 447, 154, 895, 1420, 1484, 3657, 9659, 1769, 249, 2484, 8348, 2187, 137, 6876, 292, 1635, 19026, 1213, 35, 1516, 843, 221, 721, 523, 126, 4580, 412201, 1767701)
for ($i = 0; $i -lt $blValues.Count; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 64).Value = $blValues[$i]
}

# 4. Footer note (row 266) - bump the "next cumulative update" date.
$ws.Cells.Item(266, 1).Value = "5. This file will be updated daily; the next cumulative update will be 6/24/2020."

# 5. Match the new column's width to its neighbor.
$ws.Columns.Item(64).ColumnWidth = 12

# 6. Touch UsedRange so the sheet dimension / per-row `spans` metadata is
#    recalculated to include the new column BL.
$u = $ws.UsedRange
